# Add new localization entries for "goal" related strings (volume, total
# volume, efficiency, and several error messages) to the "en" language sheet.
#
# The write order below is chosen so that new shared-string table entries
# are appended in the same order Excel produced them in the reference
# edit: the two "value" strings for the (future) rows 10/11 are entered
# first, followed by rows 7-9 in order, then row 12, and finally the two
# remaining "key" strings for rows 11 and 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "Not enough volume!"
$ws.Range("B11").Value = "Each object's height must be {0} tall!"

$ws.Range("A7").Value = "goal_volume_req"
$ws.Range("B7").Value = "Volume:"

$ws.Range("A8").Value = "goal_total_volume"
$ws.Range("B8").Value = "Total Volume:"

$ws.Range("A9").Value = "goal_efficiency"
$ws.Range("B9").Value = "Efficiency:"

$ws.Range("A12").Value = "goal_error_not_found"
$ws.Range("B12").Value = "No matching objects found!"

$ws.Range("A11").Value = "goal_error_height_not_met"
$ws.Range("A10").Value = "goal_error_volume_not_enough"

# Match the final selection/active cell recorded in the workbook (A12).
$ws.Range("A12").Select()
